# training_diary.xlsx: add support for preset & sticky values (closes #32)
#
# Adds a new "value" column (J) to the survey sheet used to mark fields
# whose answer should be preset (e.g. head(squats,1)) or simply kept
# "sticky" from the previous entry. Also refreshes the mood icons to the
# outlined/3x FontAwesome variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- New column J header: "value" -----------------------------------
$ws.Range("J1").Value = "value"

# --- J6: preset value for squats (pull the last entered squats value) -
$ws.Range("J6").Value = "head(squats,1)"

# --- Sticky flags: keep previous answer as the new default -----------
$ws.Range("J3").Value = "sticky"
$ws.Range("J5").Value = "sticky"
$ws.Range("J7").Value = "sticky"
$ws.Range("J8").Value = "sticky"

# J4 (situps) is also sticky, but left-aligned (new style) -------------
$ws.Range("J4").Value = "sticky"
$ws.Range("J4").HorizontalAlignment = -4131
$ws.Range("J4").WrapText = $true

# --- Refresh mood icons to outlined / 3x FontAwesome glyphs -----------
$ws.Range("H10").Value = '<i class="fa fa-frown-o fa-3x"></i>'
$ws.Range("G10").Value = '<i class="fa fa-smile-o fa-3x"></i>'

# Row 10 (mood icons row) shrinks now that labels changed
$ws.Rows.Item(10).RowHeight = 45

# --- Cursor / view bookkeeping, mirroring the saved selection ---------
$ws.Range("H16").Select()
